$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")
$ws.Range("D2").Value = "Tolstoy Inc"
$ws.Range("D3").Value = "Wild LLC"
$ws.Range("D4").Value = "Tolstoy Inc"
$ws.Range("D5").Select()
